$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: was blank, now "A"
$ws.Range("D2").Value = "A"

# H2: update Expected Behaviour text for VT200-0851 test (add Iconposition validations)
$ws.Range("H2").Value = "validate1
{
validate_PageTitle=Manual Compliance Ruby Specs
};
validate2
{
validate_PageTitle=Signal Ruby Test
};
validate3
{
validate_Text_Exists=VT200-0851
};
validate4
{
validate_Screenshot=VT200-0851
validate_Iconposition=signalview_xpath,left,20
validate_Iconposition=signalview_xpath,top,40
};"

# G2: update Steps text for VT200-0851 test (remove duplicated trailing wait/screenshot)
$ws.Range("G2").Value = "wait(3);
validate1;
SwitchApp(NATIVE_APP);
ClickNativeIcon(VT200_0851_home_xpath);
SwitchApp(WEBVIEW);
link_Click(signal_test_link);
validate2;
SelectTestToRun(VT200_0851_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(2);
TakeScreenshot(VT200-0851);
validate4;"

# J3: clear "Pass"
$ws.Range("J3").ClearContents()

# Row 2 height increased to fit new (longer) text
$ws.Rows(2).RowHeight = 203.25

# Update selected cell to D1
$ws.Range("D1").Select()

$wb.Save()
